$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date-header columns AD (2/16/20) and AE (2/17/20) -----------------
# Force text (not auto-parsed dates) to match the other "m/d/yy" headers,
# then give them the same bold/centered/bordered look as the rest of row 1.
$hdr = $ws.Range("AD1:AE1")
$hdr.NumberFormat = "@"
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$ws.Range("AD1").Value = "2/16/20"
$ws.Range("AE1").Value = "2/17/20"

# --- New daily confirmed-case counts for 2/16/20 and 2/17/20 ---------------
$data = @(
    @{Row=2; AD=962; AE=973},
    @{Row=3; AD=380; AE=381},
    @{Row=4; AD=551; AE=553},
    @{Row=5; AD=287; AE=290},
    @{Row=6; AD=90; AE=91},
    @{Row=7; AD=1316; AE=1322},
    @{Row=8; AD=237; AE=238},
    @{Row=9; AD=144; AE=146},
    @{Row=10; AD=162; AE=163},
    @{Row=11; AD=300; AE=301},
    @{Row=12; AD=445; AE=457},
    @{Row=13; AD=1231; AE=1246},
    @{Row=14; AD=58182; AE=59989},
    @{Row=15; AD=1004; AE=1006},
    @{Row=16; AD=70; AE=72},
    @{Row=17; AD=617; AE=626},
    @{Row=18; AD=925; AE=930},
    @{Row=19; AD=89; AE=89},
    @{Row=20; AD=121; AE=121},
    @{Row=21; AD=70; AE=70},
    @{Row=22; AD=18; AE=18},
    @{Row=23; AD=236; AE=240},
    @{Row=24; AD=537; AE=541},
    @{Row=25; AD=328; AE=333},
    @{Row=26; AD=129; AE=130},
    @{Row=27; AD=481; AE=495},
    @{Row=28; AD=124; AE=125},
    @{Row=29; AD=1; AE=1},
    @{Row=30; AD=71; AE=75},
    @{Row=31; AD=171; AE=171},
    @{Row=32; AD=1167; AE=1171},
    @{Row=33; AD=34; AE=35},
    @{Row=34; AD=59; AE=66},
    @{Row=35; AD=29; AE=30},
    @{Row=36; AD=20; AE=22},
    @{Row=37; AD=1; AE=1},
    @{Row=38; AD=2; AE=2},
    @{Row=39; AD=1; AE=1},
    @{Row=40; AD=10; AE=10},
    @{Row=41; AD=57; AE=60},
    @{Row=42; AD=75; AE=77},
    @{Row=43; AD=16; AE=16},
    @{Row=44; AD=12; AE=12},
    @{Row=45; AD=1; AE=1},
    @{Row=46; AD=22; AE=22},
    @{Row=47; AD=2; AE=2},
    @{Row=48; AD=4; AE=5},
    @{Row=49; AD=1; AE=1},
    @{Row=50; AD=1; AE=1},
    @{Row=51; AD=4; AE=4},
    @{Row=52; AD=4; AE=4},
    @{Row=53; AD=5; AE=5},
    @{Row=54; AD=1; AE=1},
    @{Row=55; AD=1; AE=1},
    @{Row=56; AD=16; AE=16},
    @{Row=57; AD=1; AE=1},
    @{Row=58; AD=9; AE=9},
    @{Row=59; AD=3; AE=3},
    @{Row=60; AD=3; AE=3},
    @{Row=61; AD=1; AE=1},
    @{Row=62; AD=3; AE=3},
    @{Row=63; AD=9; AE=9},
    @{Row=64; AD=2; AE=2},
    @{Row=65; AD=1; AE=1},
    @{Row=66; AD=2; AE=2},
    @{Row=67; AD=2; AE=2},
    @{Row=68; AD=2; AE=2},
    @{Row=69; AD=1; AE=1},
    @{Row=70; AD=2; AE=2},
    @{Row=71; AD=1; AE=1},
    @{Row=72; AD=1; AE=1},
    @{Row=73; AD=355; AE=454},
    @{Row=74; AD=2; AE=2},
    @{Row=75; AD=1; AE=1},
    @{Row=76; AD=1; AE=1}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 30).Value = $item.AD
    $ws.Cells.Item($item.Row, 31).Value = $item.AE
}
